$d = $word.ActiveDocument

function Append-Run($para, [string]$text, [bool]$bold) {
    $pr = $para.Range
    $textLen = $pr.End - $pr.Start - 1
    $insertPos = $pr.Start + $textLen
    $ip = $d.Range($insertPos, $insertPos)
    $ip.InsertAfter($text)
    $newRange = $d.Range($insertPos, $insertPos + $text.Length)
    if ($bold) {
        $newRange.Font.Bold = 1
    } else {
        $newRange.Font.Bold = 0
    }
    return $newRange
}

# ------------------------------------------------------------------
# 1. After paragraph 18 ("clear = limpar terminal"), add 5 new
#    paragraphs (git status / empty / git init / git add <file> /
#    git add .), then fill the final (formerly trailing empty)
#    paragraph with the "git commit" content + move the _GoBack
#    bookmark there.
# ------------------------------------------------------------------

$pClear = $d.Paragraphs(18)
$pClear.Range.InsertParagraphAfter()

# --- paragraph 19: git status = verificar se houve alguma alteração ou algo do tipo.
$p19 = $d.Paragraphs(19)
Append-Run $p19 "git status =" $true
Append-Run ($d.Paragraphs(19)) " verificar se houve alguma alteração ou algo do tipo. " $false

# --- paragraph 20: empty separator
$p19b = $d.Paragraphs(19)
$p19b.Range.InsertParagraphAfter()

# --- paragraph 21: git init = inicia o git ou iniciar um repositório git vazio.
$p20 = $d.Paragraphs(20)
$p20.Range.InsertParagraphAfter()
$p21 = $d.Paragraphs(21)
Append-Run $p21 "git init =" $true
Append-Run ($d.Paragraphs(21)) " inicia o git ou iniciar um repositório git vazio. " $false

# --- paragraph 22: git add <file> = adicionar um arquivo para salvar
$p21b = $d.Paragraphs(21)
$p21b.Range.InsertParagraphAfter()
$p22 = $d.Paragraphs(22)
Append-Run $p22 "git add <file> =" $true
Append-Run ($d.Paragraphs(22)) " adicionar um arquivo para salvar" $false

# --- paragraph 23: git add . = adicionar todos os arquivos para salvamento.
$p22b = $d.Paragraphs(22)
$p22b.Range.InsertParagraphAfter()
$p23 = $d.Paragraphs(23)
Append-Run $p23 "git add . =" $true
Append-Run ($d.Paragraphs(23)) " adicionar todos os arquivos para salvamento. " $false

# --- paragraph 24 is already the original trailing empty paragraph, which
#     has been pushed down to index 24 by the four InsertParagraphAfter
#     calls above - no new paragraph break needed here, just fill it in.
$p24 = $d.Paragraphs(24)
Append-Run $p24 "git commit -m “comentario” = " $true
Append-Run ($d.Paragraphs(24)) "após utilizar o git add para adicionar o arquivo, utilizamos este código para estar realizando o salvamento e colocando um comentário. " $false

# Move the _GoBack bookmark from paragraph 12 to the start of paragraph 24
# (right before "git commit"). Adding a bookmark with the same name moves it
# (Word enforces unique bookmark names), so the old one disappears automatically.
$p24fresh = $d.Paragraphs(24)
$bookmarkStartPos = $p24fresh.Range.Start
$bookmarkTextLen = "git commit -m “comentario” = ".Length
$bookmarkEndPos = $bookmarkStartPos + $bookmarkTextLen
$bookmarkRange = $d.Range($bookmarkStartPos, $bookmarkEndPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

Write-Output "DONE"
Write-Output $d.Paragraphs.Count
